# [BI-1613] Update TAF to include term type
#
# Adds a new "Term Type" column (column R) to the Template sheet of the
# trait-import test workbook, with a sample value of "germplasm passport"
# on the first data row (row 3), mirroring the other header/sample-data
# columns already present on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# New header cell, styled like the other header cells in row 1.
$ws.Range("R1").Value = "Term Type"

# New sample value on the second data row (row 3), left unstyled like the
# other plain data cells in that row.
$ws.Range("R3").Value = "germplasm passport"

# Reflect the new column in the current selection/scroll position, as a
# user editing the sheet near the new column would leave it.
[void]$ws.Range("Q7").Select()
